$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates derived from the diff (row, column letter, new value).
$updates = @(
    @{ Row=2; Col='D'; Value='28.263.19' },
    @{ Row=2; Col='E'; Value='  +0.84%  ' },
    @{ Row=3; Col='D'; Value='1.885.38' },
    @{ Row=3; Col='E'; Value='  +1.21%  ' },
    @{ Row=4; Col='D'; Value='1.008' },
    @{ Row=4; Col='E'; Value='  +0.38%  ' },
    @{ Row=5; Col='D'; Value='315.59' },
    @{ Row=5; Col='E'; Value='  +1.06%  ' },
    @{ Row=6; Col='D'; Value='1.009' },
    @{ Row=6; Col='E'; Value='  +0.65%  ' },
    @{ Row=7; Col='D'; Value='0.5150' },
    @{ Row=7; Col='E'; Value='  +1.04%  ' },
    @{ Row=8; Col='D'; Value='0.3919' },
    @{ Row=8; Col='E'; Value='  +1.86%  ' },
    @{ Row=9; Col='D'; Value='0.08410' },
    @{ Row=9; Col='E'; Value='  +1.37%  ' },
    @{ Row=10; Col='D'; Value='1.125' },
    @{ Row=10; Col='E'; Value='  +1.05%  ' },
    @{ Row=11; Col='D'; Value='41.75' },
    @{ Row=11; Col='E'; Value='  +0.52%  ' },
    @{ Row=12; Col='D'; Value='6.265' },
    @{ Row=12; Col='E'; Value='  +0.48%  ' },
    @{ Row=13; Col='B'; Value='Solana' },
    @{ Row=13; Col='C'; Value='https://coinranking.com/coin/zNZHO_Sjf+solana-sol' },
    @{ Row=13; Col='D'; Value='20.67' },
    @{ Row=13; Col='E'; Value='  +0.33%  ' },
    @{ Row=14; Col='B'; Value='WrappedEther' },
    @{ Row=14; Col='C'; Value='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Row=14; Col='D'; Value='1.879.66' },
    @{ Row=14; Col='E'; Value='  +0.92%  ' },
    @{ Row=15; Col='D'; Value='7.286' },
    @{ Row=15; Col='E'; Value='  +0.67%  ' },
    @{ Row=16; Col='D'; Value='1.005' },
    @{ Row=16; Col='E'; Value='  +0.16%  ' },
    @{ Row=17; Col='D'; Value='0.00001107' },
    @{ Row=17; Col='E'; Value='  +0.77%  ' },
    @{ Row=18; Col='D'; Value='91.42' },
    @{ Row=18; Col='E'; Value='  +0.58%  ' },
    @{ Row=19; Col='D'; Value='0.06720' },
    @{ Row=19; Col='E'; Value='  +1.31%  ' },
    @{ Row=20; Col='D'; Value='17.86' },
    @{ Row=20; Col='E'; Value='  +0.85%  ' },
    @{ Row=21; Col='D'; Value='1.009' },
    @{ Row=21; Col='E'; Value='  +0.59%  ' },
    @{ Row=22; Col='D'; Value='6.056' },
    @{ Row=22; Col='E'; Value='  +0.22%  ' },
    @{ Row=23; Col='D'; Value='28.325.24' },
    @{ Row=23; Col='E'; Value='  +0.96%  ' },
    @{ Row=24; Col='D'; Value='11.16' },
    @{ Row=24; Col='E'; Value='  +0.40%  ' },
    @{ Row=25; Col='D'; Value='2.251' },
    @{ Row=25; Col='E'; Value='  +0.81%  ' },
    @{ Row=26; Col='D'; Value='159.67' },
    @{ Row=26; Col='E'; Value='  +1.18%  ' },
    @{ Row=27; Col='D'; Value='2.473' },
    @{ Row=27; Col='E'; Value='  -2.56%  ' },
    @{ Row=28; Col='D'; Value='20.79' },
    @{ Row=28; Col='E'; Value='  +1.28%  ' },
    @{ Row=29; Col='D'; Value='126.77' },
    @{ Row=29; Col='E'; Value='  +1.46%  ' },
    @{ Row=30; Col='D'; Value='0.1057' },
    @{ Row=30; Col='E'; Value='  +0.05%  ' },
    @{ Row=31; Col='D'; Value='1.038' },
    @{ Row=31; Col='E'; Value='  +0.01%  ' },
    @{ Row=32; Col='D'; Value='5.897' },
    @{ Row=32; Col='E'; Value='  +0.11%  ' },
    @{ Row=33; Col='D'; Value='3.621' },
    @{ Row=33; Col='E'; Value='  +0.72%  ' },
    @{ Row=34; Col='D'; Value='9.612' },
    @{ Row=34; Col='E'; Value='  +2.15%  ' },
    @{ Row=35; Col='D'; Value='0.02459' },
    @{ Row=35; Col='E'; Value='  +1.48%  ' },
    @{ Row=36; Col='D'; Value='0.06581' },
    @{ Row=36; Col='E'; Value='  +0.61%  ' },
    @{ Row=37; Col='D'; Value='0.2219' },
    @{ Row=37; Col='E'; Value='  +2.05%  ' },
    @{ Row=38; Col='D'; Value='1.199' },
    @{ Row=38; Col='E'; Value='  -0.46%  ' },
    @{ Row=39; Col='D'; Value='0.6512' },
    @{ Row=39; Col='E'; Value='  +0.65%  ' },
    @{ Row=40; Col='D'; Value='1.238' },
    @{ Row=40; Col='E'; Value='  +0.92%  ' },
    @{ Row=41; Col='D'; Value='5.027' },
    @{ Row=41; Col='E'; Value='  +0.74%  ' },
    @{ Row=42; Col='D'; Value='11.33' },
    @{ Row=42; Col='E'; Value='  +1.18%  ' },
    @{ Row=43; Col='D'; Value='0.6114' },
    @{ Row=43; Col='E'; Value='  -0.23%  ' },
    @{ Row=44; Col='D'; Value='13.19' },
    @{ Row=44; Col='E'; Value='  +0.22%  ' },
    @{ Row=45; Col='D'; Value='3.693' },
    @{ Row=45; Col='E'; Value='  +1.01%  ' },
    @{ Row=46; Col='D'; Value='1.282' },
    @{ Row=46; Col='E'; Value='  -0.38%  ' },
    @{ Row=47; Col='D'; Value='2.018' },
    @{ Row=47; Col='E'; Value='  +0.13%  ' },
    @{ Row=48; Col='D'; Value='1.236' },
    @{ Row=48; Col='E'; Value='  +2.21%  ' },
    @{ Row=49; Col='D'; Value='121.56' },
    @{ Row=49; Col='E'; Value='  +1.15%  ' },
    @{ Row=50; Col='D'; Value='0.06921' },
    @{ Row=50; Col='E'; Value='  +1.08%  ' },
    @{ Row=51; Col='D'; Value='78.05' },
    @{ Row=51; Col='E'; Value='  -0.48%  ' }
)

foreach ($u in $updates) {
    $ws.Range("$($u.Col)$($u.Row)").Value = $u.Value
}
